# Darks Knight Description.docx edit script
# Applies the textual changes described by the commit:
#   "Conditional compile for 64 or 32 bit TheSky Added version to form title"

$d = $word.ActiveDocument

# Unicode helpers for curly quotes used in the source document.
$lq = [char]0x201C
$rq = [char]0x201D

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found -> $find"
    }
}

# 1. "applet uses TSX CAO to create a" -> "automation application uses The Sky to create a"
Replace-Text " applet uses TSX CAO to create a" " automation application uses The Sky to create a"

# 2. "of dark frames, and, with Rev 1.3, bias frames as well" -> "of dark and bias frames"
Replace-Text " of dark frames, and, with Rev 1.3, bias frames as well" " of dark and bias frames"

# 3. "the Take Series CAO function" -> "the Take Series function"
Replace-Text " the Take Series CAO function" " the Take Series function"

# 4. Application renamed from "Night Shift" to "Darks Knight"
Replace-Text "Night Shift is a Windows Forms executable" "Darks Knight is a Windows Forms executable"

# 5. Conditional 32/64-bit installer zip names
$findZip = "Download " + $lq + "DarkKnights.zip" + $rq + " from the"
$replZip = "Download " + $lq + "DarkKnights32.zip" + $rq + " or " + $lq + "DarkKnights64.zip" + $rq + " from the"
Replace-Text $findZip $replZip

# 6. Start menu category / shortcut name corrections
#    (kept quote-free so the existing straight quotes around the names
#    are left untouched instead of being smart-quoted by Find/Replace)
Replace-Text "TXTToolkit" "TSX Toolkit"
Replace-Text "DarksKnight" "Darks Knight"

# 7. "TSX settings" -> "THE SKY settings"
Replace-Text "according to TSX settings" "according to THE SKY settings"

# 8. "other TSX applets" -> "other THE SKY applets"
Replace-Text "used for other TSX applets" "used for other THE SKY applets"

# 9. Remove stale _GoBack bookmark left over from the prior save.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
